# Re-case the schema header row labels from snake_case to the schematic
# PascalCase/ID-style casing, add the new "HasParticipant" column to the
# Study sheets, and fix the "race" option casing in the data-validation
# list for the Participant sheet (CLI + schema-compliance casing pass).

$wb = $excel.ActiveWorkbook

# --- 1. Header row relabels (applies identically to each sheet and its
#        duplicate "<Name>1" counterpart) ------------------------------

$headerMap = @{
    "Assay"  = @("UsesBiospecimen", "HasOutput");
    "Assay1" = @("UsesBiospecimen", "HasOutput");

    "Biospecimen"  = @(
        "AgeAtBiospecimenCollection", "BiospecimenStorage", "CollectionID",
        "CollectionSampleType", "ContainerID", "HasDatafile", "HasParticipant",
        "HasStudy", "LaboratoryProcedure", "ParentSampleID", "ParentSampleType",
        "SampleAvailability", "SampleID", "SampleType", "Volume", "VolumeUnit",
        "HasAliquot"
    );
    "Biospecimen1" = @(
        "AgeAtBiospecimenCollection", "BiospecimenStorage", "CollectionID",
        "CollectionSampleType", "ContainerID", "HasDatafile", "HasParticipant",
        "HasStudy", "LaboratoryProcedure", "ParentSampleID", "ParentSampleType",
        "SampleAvailability", "SampleID", "SampleType", "Volume", "VolumeUnit",
        "HasAliquot"
    );

    "Condition"  = @(
        "HasParticipant", "AgeAtConditionObservation", "MONDOLabel", "MONDOCode",
        "ConditionInterpretation", "ConditionDataSource", "HPOLabel", "HPOCode",
        "MAXOLabel", "MAXOCode", "OtherLabel", "OtherCode"
    );
    "Condition1" = @(
        "HasParticipant", "AgeAtConditionObservation", "MONDOLabel", "MONDOCode",
        "ConditionInterpretation", "ConditionDataSource", "HPOLabel", "HPOCode",
        "MAXOLabel", "MAXOCode", "OtherLabel", "OtherCode"
    );

    "DataFile"  = @(
        "AccessURL", "CollectionID", "DataAccess", "DataCategory", "DataType",
        "ExperimentalStrategy", "FileID", "FileName", "Format", "HasBiospecimen",
        "HasParticipant", "HasStudy", "ParticipantID", "Size", "OriginalFileName"
    );
    "DataFile1" = @(
        "AccessURL", "CollectionID", "DataAccess", "DataCategory", "DataType",
        "ExperimentalStrategy", "FileID", "FileName", "Format", "HasBiospecimen",
        "HasParticipant", "HasStudy", "ParticipantID", "Size", "OriginalFileName"
    );

    "FamilyGroup"  = @("HasParticipant");
    "FamilyGroup1" = @("HasParticipant");

    "Participant"  = @(
        "AgeAtLastVitalStatus", "DownSyndromeStatus", "Ethnicity", "ExternalID",
        "FamilyID", "FamilyRelationship", "FamilyType", "FatherID", "HasDatafile",
        "HasStudy", "MotherID", "OutcomesVitalStatus", "ParticipantID", "Race", "Sex"
    );
    "Participant1" = @(
        "AgeAtLastVitalStatus", "DownSyndromeStatus", "Ethnicity", "ExternalID",
        "FamilyID", "FamilyRelationship", "FamilyType", "FatherID", "HasDatafile",
        "HasStudy", "MotherID", "OutcomesVitalStatus", "ParticipantID", "Race", "Sex"
    );

    "Study"  = @("dbGap", "Program", "StudyCode", "StudyName");
    "Study1" = @("dbGap", "Program", "StudyCode", "StudyName");
}

foreach ($sheetName in $headerMap.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $values = $headerMap[$sheetName]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item(1, $i + 1).Value = $values[$i]
    }
}

# --- 2. New "HasParticipant" column on the Study sheets (A1:D1 -> A1:E1) ---

foreach ($sheetName in @("Study", "Study1")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Cells.Item(1, 5).Value = "HasParticipant"
}

# --- 3. Data validation list fix: more_than_one_race -> more_than_one_Race
#        on the Participant sheet's "Race" column (both duplicated
#        dataValidation entries for N2:N1048576 share the same formula) ---

$raceFormula = '"american_indian_or_alaskan_native,asian,black_or_african_american,more_than_one_Race,native_hawaiian_or_other_pacific_islander,other,white,prefer_not_to_answer"'

foreach ($sheetName in @("Participant", "Participant1")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rng = $ws.Range("N2:N1048576")
    if ($rng.Validation.Formula1 -ne $null) {
        $rng.Validation.Delete()
        $rng.Validation.Add(3, 1, 1, $raceFormula)
    }
}
